# Applies the weekly work-report update:
#  - refresh "Report Generated On" timestamp
#  - refresh Total Billed Amount / Total Line Items summary figures
#  - append 7 newly-completed line items just above the TOTAL row
#  - roll the TOTAL row's grand total forward to include the new lines

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header / summary refresh
# ---------------------------------------------------------------------
$ws.Range("D5").Value2 = "Report Generated On: 08/16/2025 12:47 AM"
$ws.Range("C8").Value2 = 27981.8    # Total Billed Amount
$ws.Range("C9").Value2 = 200        # Total Line Items

# ---------------------------------------------------------------------
# 2. Make room for 7 new line-item rows just above the existing TOTAL
#    row (row 214). This pushes TOTAL down to row 221 and shifts the
#    "A214:G214" TOTAL merge down to "A221:G221" automatically.
# ---------------------------------------------------------------------
$ws.Rows("214:220").Insert() | Out-Null

# Carry over the existing alternating row styling (odd template = row
# 212's style set, even template = row 213's style set) so the new
# rows visually match the rest of the table instead of picking up
# default formatting.
$ws.Range("A212:H212").Copy() | Out-Null
$ws.Range("A214:H214").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A216:H216").PasteSpecial(-4122) | Out-Null
$ws.Range("A218:H218").PasteSpecial(-4122) | Out-Null
$ws.Range("A220:H220").PasteSpecial(-4122) | Out-Null

$ws.Range("A213:H213").Copy() | Out-Null
$ws.Range("A215:H215").PasteSpecial(-4122) | Out-Null
$ws.Range("A217:H217").PasteSpecial(-4122) | Out-Null
$ws.Range("A219:H219").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Populate the 7 new completed line items
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=214; Point="Point 05"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=1; Amount=195.83 },
    @{ Row=215; Point="Point 07"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=1; Amount=195.83 },
    @{ Row=216; Point="Point 09"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=1; Amount=195.83 },
    @{ Row=217; Point="Point 11"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=1; Amount=195.83 },
    @{ Row=218; Point="Point 13"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=1; Amount=195.83 },
    @{ Row=219; Point="Point 29"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=3; Amount=587.49 },
    @{ Row=220; Point="Point 03"; Code="SVC-40-TP-ALA-PC"; Type="Trans"; Desc="SVC,4/0,Trip,Alum/Aly Neut,Paral Comm"; Unit="EA"; Qty=1; Amount=195.83 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.Point    # A: Point
    $ws.Cells.Item($row, 2).Value2 = $r.Code     # B: Scope/unit code
    $ws.Cells.Item($row, 3).Value2 = $r.Type     # C: Trans/Inst
    $ws.Cells.Item($row, 4).Value2 = $r.Desc     # D: Description
    $ws.Cells.Item($row, 5).Value2 = $r.Unit     # E: Unit of measure
    $ws.Cells.Item($row, 6).Value2 = $r.Qty      # F: Quantity
    $ws.Cells.Item($row, 8).Value2 = $r.Amount   # H: Extended amount
}

# ---------------------------------------------------------------------
# 4. Roll the TOTAL row forward (now at row 221) to the new grand total
# ---------------------------------------------------------------------
$ws.Range("H221").Value2 = 15519.15

"Update complete: TOTAL now at row 221 with 7 new line items."
